$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-24 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-01-25 Saturday", 2)

$d.Content.Find.Execute("45×36=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "35×99=", 2)

$d.Content.Find.Execute("79×25=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "41×27=", 2)

$d.Content.Find.Execute("57×60=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "94×35=", 2)

$d.Content.Find.Execute("99×69=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "68×54=", 2)

$d.Content.Find.Execute("79×16=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "93×96=", 2)

$d.Content.Find.Execute("47×36=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "87×79=", 2)

$d.Content.Find.Execute("43×93=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "88×33=", 2)

$d.Content.Find.Execute("14×13=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "35×77=", 2)

$d.Content.Find.Execute("63×54=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "90×33=", 2)

$d.Content.Find.Execute("14×82=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "26×89=", 2)

$d.Content.Find.Execute("61×90=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "80×66=", 2)

$d.Content.Find.Execute("85×67=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "13×65=", 2)

$d.Content.Find.Execute("80×30=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "66×46=", 2)

$d.Content.Find.Execute("41×97=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "41×82=", 2)

$d.Content.Find.Execute("79×66=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "74×54=", 2)

$d.Content.Find.Execute("17×97=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "86×74=", 2)

$d.Content.Find.Execute("64×52=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "69×46=", 2)

$d.Content.Find.Execute("59×83=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "65×42=", 2)

$d.Content.Find.Execute("84×26=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25×33=", 2)

$d.Content.Find.Execute("20×53=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "30×78=", 2)

$d.Content.Find.Execute("62×26=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "32×31=", 2)

$d.Content.Find.Execute("94×65=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "81×40=", 2)

$d.Content.Find.Execute("80×51=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "69×64=", 2)

$d.Content.Find.Execute("75×13=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "29×71=", 2)

$d.Content.Find.Execute("99×30=", $true, $false, $false, $false, $false,
                         $true, 1, $false, "76×66=", 2)
